$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the specific "X of a Kind in a Deck of Cards" list item that is
# immediately preceded by "...Greatest Common Divisor of Strings (solved)"
# (there are two "X of a Kind in a Deck of Cards" bullets in this document;
# only this one needs the " (solved)" annotation added).
$targetPara = $null
$sourcePara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "X of a Kind in a Deck of Cards") {
        $prev = $p.Previous()
        $prevText = $prev.Range.Text.TrimEnd([char]13, [char]7)
        if ($prevText -eq "Greatest Common Divisor of Strings (solved)") {
            $targetPara = $p
            $sourcePara = $prev
            break
        }
    }
}

if ($targetPara -ne $null) {
    # Grab the already-formatted " (solved)" run from the preceding paragraph
    # so the newly-inserted text picks up identical run formatting (Arial,
    # sz 22, auto color, etc.) instead of plain/default formatting.
    $srcRange = $sourcePara.Range
    $srcEnd = $srcRange.End - 1
    $srcStart = $srcEnd - 9
    $srcMarker = $d.Range($srcStart, $srcEnd)

    # Insert at the end of the target paragraph (just before its paragraph mark).
    $tgtRange = $targetPara.Range
    $insertPoint = $tgtRange.End - 1
    $tgt = $d.Range($insertPoint, $insertPoint)
    $tgt.FormattedText = $srcMarker.FormattedText
}
